$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.871.69"
$ws.Range("E2").Value = "  +1.40%  "

# Row 3
$ws.Range("D3").Value = "1.905.29"
$ws.Range("E3").Value = "  +2.44%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.76"
$ws.Range("E5").Value = "  +1.73%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  -0.22%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4842"
$ws.Range("E7").Value = "  +1.59%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3791"
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07366"
$ws.Range("E9").Value = "  +0.64%  "

# Row 10
$ws.Range("E10").Value = "  +0.15%  "

# Row 11
$ws.Range("E11").Value = "  +0.29%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07750"
$ws.Range("E12").Value = "  -0.47%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.888.45"
$ws.Range("E13").Value = "  +1.83%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.488"
$ws.Range("E14").Value = "  +0.72%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.636"
$ws.Range("E15").Value = "  +1.17%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.90"
$ws.Range("E16").Value = "  +1.91%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("E17").Value = "  -0.05%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008854"
$ws.Range("E18").Value = "  +0.36%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  -0.26%  "

# Row 20
$ws.Range("D20").Value = "27.926.93"
$ws.Range("E20").Value = "  +1.95%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.66"
$ws.Range("E21").Value = "  +0.20%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.166"
$ws.Range("E22").Value = "  +1.52%  "

# Row 23
$ws.Range("D23").Value = "2.204.76"
$ws.Range("E23").Value = "  +5.89%  "

# Row 24
$ws.Range("E24").Value = "  +1.82%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.921"
$ws.Range("E25").Value = "  -0.94%  "

# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.21"
$ws.Range("E26").Value = "  +0.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.47"
$ws.Range("E27").Value = "  +0.00%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.124"
$ws.Range("E28").Value = "  +5.93%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.28"
$ws.Range("E29").Value = "  +1.71%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.959"
$ws.Range("E30").Value = "  +0.24%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08966"
$ws.Range("E31").Value = "  +1.08%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.236"
$ws.Range("E32").Value = "  -2.80%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.253"
$ws.Range("E33").Value = "  +4.17%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7676"
$ws.Range("E34").Value = "  +2.09%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.650"
$ws.Range("E35").Value = "  +1.55%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02046"
$ws.Range("E36").Value = "  +0.04%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.537"
$ws.Range("E37").Value = "  -6.18%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05277"
$ws.Range("E39").Value = "  -0.33%  "

# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.999"
$ws.Range("E40").Value = "  +0.63%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5473"
$ws.Range("E41").Value = "  -1.31%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.963"
$ws.Range("E42").Value = "  -0.85%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1526"
$ws.Range("E43").Value = "  +0.45%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.439"
$ws.Range("E44").Value = "  -1.34%  "

# Row 45
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "109.25"
$ws.Range("E45").Value = "  +5.19%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.66"
$ws.Range("E46").Value = "  -0.05%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4807"
$ws.Range("E47").Value = "  -1.19%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("E48").Value = "  -0.22%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.652"
$ws.Range("E49").Value = "  -0.78%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.62"
$ws.Range("E50").Value = "  +0.60%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06091"
$ws.Range("E51").Value = "  -0.07%  "
